$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 31: Validation of an irritability measure in preschoolers...
$ws.Cells.Item(31, 1).Value = 'Validation of an irritability measure in preschoolers in school-based and clinical Brazilian samples'
$ws.Cells.Item(31, 2).Value = 'Luisa Shiguemi Sugaya, Katharina Kircanski, Argyris Stringaris, Guilherme V. Polanczyk & Ellen Leibenluft'
$ws.Cells.Item(31, 3).Value = 'European Child & Adolescent Psychiatry'
$ws.Cells.Item(31, 4).Value = 2021
$ws.Cells.Item(31, 5).Value = 'https://doi.org/10.1007/s00787-020-01701-6'
$ws.Cells.Item(31, 6).Value = 'The Affective Reactivity Index (ARI) is an irritability measure with good psychometric properties. However, there are no published studies in preschool children, an important population in which to differentiate normative from non-normative irritability. The goal of this study was to validate the ARI in preschoolers. Two samples were included: a school-based sample (N = 487, mean age = 57.80 ± 7.23 months, 52.8% male) and a clinical sample of children with Attention Deficit Hyperactivity Disorder (ADHD; N = 153, mean age = 60.5 ± 7.6 months, 83.7% males). Confirmatory factor analysis assessed ARI unidimensionality. ARI criterion validity was tested through comparison to other scales measuring irritability, related constructs, and other aspects of psychopathology. Test–retest reliability was assessed in the school-based sample. Analyses confirmed a single-factor structure and good internal consistency. The ARI showed stronger correlations with irritability measures than with measures of other constructs. In the clinical sample, ADHD children with comorbid disruptive behavior disorders had higher ARI scores than those without this comorbidity. In the school-based sample, test–retest reliability was moderate. This is the first study to demonstrate ARI validity and reliability in preschoolers. The scale performed well in both school-based and clinical samples. Having a concise and validated irritability measure for preschoolers may facilitate both clinical assessment and research on early irritability.'
$ws.Cells.Item(31, 7).Value = 'ARI evaluated'
$ws.Cells.Item(31, 8).Value = 'Psychometric'

# Row 32: Assessing beliefs about emotion generation and change (CMBQ)...
$ws.Cells.Item(32, 1).Value = 'Assessing beliefs about emotion generation and change: The conceptualisation, development, and validation of the Cognitive Mediation Beliefs Questionnaire (CMBQ)'
$ws.Cells.Item(32, 2).Value = 'Martin J. Turner, Andrew G. Wood, Daniel Boatwright, Nanaki Chadha, Jennifer K. Jones, Richard Bennett'
$ws.Cells.Item(32, 3).Value = 'Psychotherapy Research'
$ws.Cells.Item(32, 4).Value = 2021
$ws.Cells.Item(32, 5).Value = 'https://doi.org/10.1080/10503307.2020.1871524'
$ws.Cells.Item(32, 6).Value = 'The ability to regulate emotions is important for human function and health. That emotion regulation can be achieved through cognitive change is predicated on the notion of cognitive mediation. However, the extent to which individuals believe that their emotions are cognitively mediated (C–M), or in contrast, that their emotions occur via stimulus-response (S-R), is underexplored, and whether C–M and S-R beliefs shape emotion reactivity is not yet known. Research that addresses these empirical needs could inform emotion regulation interventions such as cognitive behavioural therapies (CBTs). The current paper reports the development and initial validity testing of the cognitive mediation beliefs questionnaire (CMBQ). Five studies report the factor structure, the construct and criterion validity, and the test-retest reliability of the CMBQ. The CMBQ was found to have a correlated two-factor structure (C–M change beliefs, and S-R generation beliefs). Higher C–M change beliefs and lower S-R generation beliefs were related to greater emotion regulation, greater thought control ability, higher positive mental health, and lower emotion reactivity. The CMBQ also demonstrated acceptable test-retest reliability. Initial testing indicates that the CMBQ is a valid and reliable questionnaire for psychometric use in adult populations, including those with a diagnosed mental health condition.'
$ws.Cells.Item(32, 7).Value = 'ARI referenced in paper'
$ws.Cells.Item(32, 8).Value = 'Psychometric'

